$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest crypto data.
# Price values are forced to text via NumberFormat "@" (then restored to
# the "Normal" style) so Excel does not reinterpret them as numbers and
# lose their original formatting (e.g. thousand separators, trailing zeros).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.572.06"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.03%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.248.26"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -1.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.134"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.421"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.812.83"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.57%  "
$ws.Range("E13").Value = "  +0.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "68.536.87"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.04%  "
$ws.Range("E16").Value = "  +1.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.256.76"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("E18").Value = "  -0.68%  "
$ws.Range("E19").Value = "  -1.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "393.77"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.70%  "
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.36"
$ws.Range("D23").Style = "Normal"
$ws.Range("E24").Value = "  +0.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000118"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.67%  "
$ws.Range("E26").Value = "  +4.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.60"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.77%  "
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.67"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.93"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.19%  "
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.27"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "164.30"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.56%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.49"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.90%  "
$ws.Range("E37").Value = "  +4.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.821"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "26.25"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.55"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "41.42"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.15%  "
$ws.Range("E43").Value = "  -6.21%  "
$ws.Range("E44").Value = "  +0.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "343.18"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.590.40"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.64"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.48%  "
$ws.Range("E48").Value = "  -0.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "31.91"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.51%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.89%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.101"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.33%  "
